$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix "Objetivos:" row (row 10) content: the text that used to be the long
# "Apresentar um quadro conceitual..." paragraph becomes "5701460 - Antonio Iacono"
$ws.Range("B10").Value = "5701460 - Antonio Iacono"
$ws.Range("C10").Value = "5701460 - Antonio Iacono"

# --- Rebuild rows 13-24 (old) into the new rows 13-23 layout.
# Delete the old rows entirely first (also removes the now-unreferenced shared
# strings for "5701460 - Antonio Iacono" duplicate, "Caracterização...", the
# old "Critério:"/"Norma de recuperação:"/"Bibliografia:" bodies, etc.)
$ws.Range("A13:C24").EntireRow.Delete()

# Row 13: Programa resumido: / Semestral / Semestral (ht 60)
$ws.Cells.Item(13,1).Value = "Programa resumido:"
$ws.Cells.Item(13,2).Value = "Semestral"
$ws.Cells.Item(13,3).Value = "Semestral"
$ws.Rows(13).RowHeight = 60

# Row 14: Short syllabus: / short-syllabus text (ht 60)
$ws.Cells.Item(14,1).Value = "Short syllabus:"
$ws.Cells.Item(14,2).Value = "1. Characterization of production programming and control. 2. Demand management. 3. Aggregate Production Planning. 4. Master Production Schedulling. 5. Inventory planning and control. 6. Material Requirement Planning (MRP). 7. Production Schedulling. 8. Detailed scheduling of production. 9. Just In Time (JIT). 10. Theory of Constraints (TOC). 11. Production control systems."
$ws.Cells.Item(14,3).Value = "1. Characterization of production programming and control. 2. Demand management. 3. Aggregate Production Planning. 4. Master Production Schedulling. 5. Inventory planning and control. 6. Material Requirement Planning (MRP). 7. Production Schedulling. 8. Detailed scheduling of production. 9. Just In Time (JIT). 10. Theory of Constraints (TOC). 11. Production control systems."
$ws.Rows(14).RowHeight = 60

# Row 15: Programa: / 01/01/2021 / 01/01/2021 (ht 120)
# (use a value-only paste from an existing "01/01/2021" cell so the date-like
# text is not auto-converted into a date serial number)
$ws.Cells.Item(15,1).Value = "Programa:"
$ws.Rows(15).RowHeight = 120
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Row 16: Syllabus: / short-syllabus text (ht 120)
$ws.Cells.Item(16,1).Value = "Syllabus:"
$ws.Cells.Item(16,2).Value = "1. Characterization of production programming and control. 2. Demand management. 3. Aggregate Production Planning. 4. Master Production Schedulling. 5. Inventory planning and control. 6. Material Requirement Planning (MRP). 7. Production Schedulling. 8. Detailed scheduling of production. 9. Just In Time (JIT). 10. Theory of Constraints (TOC). 11. Production control systems."
$ws.Cells.Item(16,3).Value = "1. Characterization of production programming and control. 2. Demand management. 3. Aggregate Production Planning. 4. Master Production Schedulling. 5. Inventory planning and control. 6. Material Requirement Planning (MRP). 7. Production Schedulling. 8. Detailed scheduling of production. 9. Just In Time (JIT). 10. Theory of Constraints (TOC). 11. Production control systems."
$ws.Rows(16).RowHeight = 120

# Row 17: Avaliação: only (no height override)
$ws.Cells.Item(17,1).Value = "Avaliação:"

# Row 18: Método: / 5701460 - Antonio Iacono / 5701460 - Antonio Iacono (ht 60)
$ws.Cells.Item(18,1).Value = "Método:"
$ws.Cells.Item(18,2).Value = "5701460 - Antonio Iacono"
$ws.Cells.Item(18,3).Value = "5701460 - Antonio Iacono"
$ws.Rows(18).RowHeight = 60

# Row 19: Critério: / Aulas expositivas... (ht 60)
$ws.Cells.Item(19,1).Value = "Critério:"
$ws.Cells.Item(19,2).Value = "Aulas expositivas teóricas, aulas práticas, aulas de exercícios. MANTIDO"
$ws.Cells.Item(19,3).Value = "Aulas expositivas teóricas, aulas práticas, aulas de exercícios. MANTIDO"
$ws.Rows(19).RowHeight = 60

# Row 20: Norma de recuperação: / M = (0,8P ... (ht 60)
$ws.Cells.Item(20,1).Value = "Norma de recuperação:"
$ws.Cells.Item(20,2).Value = "M = (0,8P + 0,2T)P = média aritmética de duas provas escritasT = Média das notas de trabalhos e exercíciosM = Média de aproveitamento do alunoAprovação com média de aproveitamento maior ou igual a 5,0 e no mínimo 70% de frequência às aulas."
$ws.Cells.Item(20,3).Value = "M = (0,8P + 0,2T)P = média aritmética de duas provas escritasT = Média das notas de trabalhos e exercíciosM = Média de aproveitamento do alunoAprovação com média de aproveitamento maior ou igual a 5,0 e no mínimo 70% de frequência às aulas."
$ws.Rows(20).RowHeight = 60

# Row 21: Bibliografia: / MF = (0,5 M ... (ht 120)
$ws.Cells.Item(21,1).Value = "Bibliografia:"
$ws.Cells.Item(21,2).Value = "MF = (0,5 M + 0,5 R)M = Média de aproveitamento do aluno, antes da recuperaçãoR = Nota de uma prova de recuperaçãoMF = nota final de aproveitamento, após a recuperaçãoAprovação com média final de aproveitamento maior ou igual a 5,0.A recuperação deverá consistir de uma prova escrita englobando a matéria toda do semestre.Terá direito à prova de recuperação aqueles alunos reprovados com nota acima de 3,0 e frequência mínima de 70%."
$ws.Cells.Item(21,3).Value = "MF = (0,5 M + 0,5 R)M = Média de aproveitamento do aluno, antes da recuperaçãoR = Nota de uma prova de recuperaçãoMF = nota final de aproveitamento, após a recuperaçãoAprovação com média final de aproveitamento maior ou igual a 5,0.A recuperação deverá consistir de uma prova escrita englobando a matéria toda do semestre.Terá direito à prova de recuperação aqueles alunos reprovados com nota acima de 3,0 e frequência mínima de 70%."
$ws.Rows(21).RowHeight = 120

# Row 22: Requisitos: only (no height override)
$ws.Cells.Item(22,1).Value = "Requisitos:"

# Row 23: LOQ4205 requirement text in B/C only (ht 30)
$ws.Cells.Item(23,2).Value = "LOQ4205 -  Sistemas Produtivos II  (Requisito fraco)`n"
$ws.Cells.Item(23,3).Value = "LOQ4205 -  Sistemas Produtivos II  (Requisito fraco)`n"
$ws.Rows(23).RowHeight = 30

# --- Fix column-B styles: newly-created cells in column B default to the
# column-A style instead of column B's style, so repair via PasteSpecial of
# formats copied from an existing, correctly-styled column B cell. Only the
# rows that actually carry a B-column value need this (17 and 22 have none).
$ws.Range("B2").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("B23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
